$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("F3").Value = "ONGOING"
$ws.Range("G4").Select()
